$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 71): N / Texp 1 / Texp 2 ------------------------
$ws.Range("B71").Value = "N"
$ws.Range("C71").Value = "Texp 1"
$ws.Range("D71").Value = "Texp 2"

# --- Bring over the number formatting used by the analogous block above --
# (B51:E59 mirrors the shape/format of the new B72:E80 block: integer N in
# column B, 0.000-formatted values in C/D, and a 0.00000-formatted (empty)
# column E).
$ws.Range("B51:E59").Copy() | Out-Null
$ws.Range("B72").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Data rows 72-80 -------------------------------------------------------
$data = @(
    @(72, 100000,  0.001,  0),
    @(73, 200000,  0.001,  0),
    @(74, 400000,  0.002,  0),
    @(75, 600000,  0.006,  0),
    @(76, 800000,  0.006,  0),
    @(77, 1200000, 0.015,  0),
    @(78, 1400000, 0.016,  0.015),
    @(79, 1600000, 0.012,  0.015),
    @(80, 1800000, 0.016,  0.015)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
}

# --- Selection left where the author left it off --------------------------
$ws.Range("C78").Select() | Out-Null
